$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (Oyuncu Adı / Pozisyon / Takım) replacing the previous table contents.
$data = @(
    @("Jaden McDaniels",     "SF,PF",    "Minnesota Timberwolves"),
    @("Malik Monk",          "PG,SG,SF", "Sacramento Kings"),
    @("Ryan Dunn",           "SF",       "Phoenix Suns"),
    @("Bam Adebayo",         "C",        "Miami Heat"),
    @("Julius Randle",       "PF,C",     "Minnesota Timberwolves"),
    @("Damian Lillard",      "PG",       "Milwaukee Bucks"),
    @("Cade Cunningham",     "PG,SG",    "Detroit Pistons"),
    @("Brandon Miller",      "SG,SF,PF", "Charlotte Hornets"),
    @("Yves Missi",          "C",        "New Orleans Pelicans"),
    @("Anthony Davis",       "PF,C",     "Los Angeles Lakers"),
    @("Herbert Jones",       "SF,PF",    "New Orleans Pelicans"),
    @("Derrick White",       "PG,SG",    "Boston Celtics"),
    @("LaMelo Ball",         "PG,SG",    "Charlotte Hornets"),
    @("Collin Sexton",       "PG,SG",    "Utah Jazz"),
    @("Isaiah Hartenstein",  "C",        "Oklahoma City Thunder"),
    @("Cameron Johnson",     "SF,PF",    "Brooklyn Nets")
)

# Clear out the previous table body (old data went through row 19) so no stale rows remain.
$ws.Range("A2:C19").Clear()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
